# Reorder the data rows (2-16) of the resources_info sheet according to the
# new ordering while keeping each row's A/B/C/D values intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 16

# Snapshot the current rows keyed by the value in column A so we can
# reassemble them in the desired order afterwards.
$rowsByKey = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $key = $ws.Cells.Item($r, 1).Value2
    $rowsByKey[$key] = @(
        $ws.Cells.Item($r, 1).Value2,
        $ws.Cells.Item($r, 2).Value2,
        $ws.Cells.Item($r, 3).Value2,
        $ws.Cells.Item($r, 4).Value2
    )
}

# Desired new order of resource keys (column A).
$newOrder = @(
    "hbw_2",
    "hw_1",
    "sm_2",
    "dm_2",
    "ov_1",
    "ov_2",
    "pm_1",
    "wt_2",
    "mm_2",
    "mm_1",
    "vgr_1",
    "wt_1",
    "sm_1",
    "hbw_1",
    "vgr_2"
)

$r = $firstRow
foreach ($key in $newOrder) {
    $data = $rowsByKey[$key]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $r++
}
